$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix country label order: swap "Namibia" and "San Vicente y las Granadinas" ---
# Row 194 currently holds Namibia's data but should be labeled
# "San Vicente y las Granadinas"; row 195 currently holds
# "San Vicente y las Granadinas"'s data (identical values) but should be
# labeled "Namibia".
$ws.Cells.Item(194, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(195, 1).Value = "Namibia"

# --- Update statistics: Brasil (row 12) ---
$ws.Cells.Item(12, 4).Value = 42991
$ws.Cells.Item(12, 5).Value = 51131

# --- Update statistics: Canada (row 15) ---
$ws.Cells.Item(15, 2).Value = 59474
$ws.Cells.Item(15, 3).Value = 2760
$ws.Cells.Item(15, 4).Value = 24908
$ws.Cells.Item(15, 5).Value = 30884
$ws.Cells.Item(15, 7).Value = 116
$ws.Cells.Item(15, 8).Value = 3682

# --- Update statistics: Costa de Marfil (row 88) ---
$ws.Cells.Item(88, 2).Value = 1398
$ws.Cells.Item(88, 3).Value = 36
$ws.Cells.Item(88, 4).Value = 653
$ws.Cells.Item(88, 5).Value = 728
$ws.Cells.Item(88, 7).Value = 2
$ws.Cells.Item(88, 8).Value = 17

# --- Update statistics: Uganda (row 158) ---
$ws.Cells.Item(158, 2).Value = 89
$ws.Cells.Item(158, 3).Value = 1
$ws.Cells.Item(158, 5).Value = 37
